$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = 4
    3 = 9
    4 = 4
    5 = 8
    6 = 6
    7 = 9
    8 = 6
    9 = 7
    10 = 4
    11 = 7
    12 = 8
    13 = 6
    14 = 4
    15 = 10
    16 = 7
    17 = 6
    18 = 4
    19 = 2
    20 = 6
    21 = 10
    22 = 3
    23 = 9
    24 = 2
    25 = 5
    26 = 10
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
